{"js": "const replacements = [\n  [\"2023-03-31 Friday\", \"2023-04-01 Saturday\"],\n  [\"30-2=\", \"89-71=\"],\n  [\"51+8=\", \"37-32=\"],\n  [\"19+39=\", \"95-12=\"],\n  [\"21-18=\", \"1+83=\"],\n  [\"68+13=\", \"91-72=\"],\n  [\"19+35=\", \"63-32=\"],\n  [\"39+19=\", \"38+27=\"],\n  [\"67-38=\", \"25-19=\"],\n  [\"99-19=\", \"47+51=\"],\n  [\"83-20=\", \"88+8=\"],\n  [\"60+7=\", \"31+25=\"],\n  [\"89-79=\", \"23+43=\"],\n  [\"62+34=\", \"0+39=\"],\n  [\"18+51=\", \"67+9=\"],\n  [\"23+35=\", \"68-0=\"],\n  [\"53-10=\", \"1+50=\"],\n  [\"54-40=\", \"59-41=\"],\n  [\"87+11=\", \"76-47=\"],\n  [\"26+59=\", \"43-22=\"],\n  [\"72-25=\", \"83-10=\"],\n  [\"82-29=\", \"86-74=\"],\n  [\"10+52=\", \"84-67=\"],\n  [\"47+0=\", \"90-41=\"],\n  [\"20+48=\", \"36+44=\"],\n  [\"33-14=\", \"23+62=\"],\n  [\"23-22=\", \"46-14=\"],\n  [\"31-21=\", \"78-18=\"],\n  [\"37-36=\", \"23+74=\"],\n  [\"10+89=\", \"81-39=\"],\n  [\"69+26=\", \"25+74=\"],\n  [\"78-52=\", \"44-6=\"],\n  [\"51+11=\", \"98-68=\"],\n  [\"98-59=\", \"70+15=\"],\n  [\"80-31=\", \"56-20=\"],\n  [\"98+1=\", \"61-28=\"],\n  [\"62-27=\", \"53+37=\"],\n  [\"67-26=\", \"37-12=\"],\n  [\"37-22=\", \"7+45=\"],\n  [\"68+23=\", \"34+23=\"],\n  [\"9+44=\", \"25+47=\"],\n  [\"68+24=\", \"25+72=\"],\n  [\"26+43=\", \"22+60=\"],\n  [\"0+88=\", \"47-1=\"],\n  [\"12+58=\", \"80-44=\"],\n  [\"75+24=\", \"31+12=\"],\n  [\"61-40=\", \"43+55=\"],\n  [\"44+3=\", \"66-25=\"],\n  [\"72-6=\", \"41-15=\"],\n  [\"25+0=\", \"21-12=\"],\n  [\"36-14=\", \"27-21=\"],\n  [\"53-24=\", \"41+5=\"],\n  [\"67-41=\", \"37-12=\"],\n  [\"75-75=\", \"85-83=\"],\n  [\"22+14=\", \"3+10=\"],\n  [\"19+49=\", \"81-19=\"],\n  [\"75-11=\", \"59+10=\"],\n  [\"16+30=\", \"95-19=\"],\n  [\"1+32=\", \"79+18=\"],\n  [\"30+34=\", \"11+5=\"],\n  [\"78-16=\", \"6+46=\"],\n  [\"14-4=\", \"21+13=\"],\n  [\"33-5=\", \"67-56=\"],\n  [\"66-6=\", \"0+65=\"],\n  [\"3+68=\", \"1+90=\"],\n  [\"66+0=\", \"77-15=\"],\n  [\"25+55=\", \"67-18=\"],\n  [\"49+27=\", \"77-14=\"],\n  [\"30+11=\", \"95-6=\"],\n  [\"8+91=\", \"79-50=\"],\n  [\"65-5=\", \"30+48=\"],\n  [\"58-1=\", \"60+39=\"],\n  [\"72-42=\", \"62+32=\"],\n  [\"43-30=\", \"97-54=\"],\n  [\"49-4=\", \"56-1=\"],\n  [\"99-97=\", \"60+6=\"],\n  [\"33+63=\", \"9+40=\"],\n  [\"74+23=\", \"96-90=\"],\n  [\"18+53=\", \"76+3=\"],\n  [\"57+21=\", \"87-33=\"],\n  [\"84-65=\", \"90-73=\"],\n  [\"44-43=\", \"36-22=\"],\n  [\"29-14=\", \"27-5=\"],\n  [\"3+23=\", \"1+58=\"],\n  [\"4+32=\", \"29+24=\"],\n  [\"23+27=\", \"54-27=\"],\n  [\"59-9=\", \"43+45=\"],\n  [\"83-42=\", \"91-63=\"],\n  [\"73-66=\", \"76-62=\"],\n  [\"0+64=\", \"62-0=\"],\n  [\"72-3=\", \"96-70=\"],\n  [\"50+39=\", \"30+30=\"],\n  [\"46+41=\", \"10+62=\"],\n  [\"64-61=\", \"70-39=\"],\n  [\"38-3=\", \"9+10=\"],\n  [\"68-32=\", \"91-58=\"],\n  [\"54-35=\", \"11+12=\"],\n  [\"19-11=\", \"27-10=\"],\n  [\"14+58=\", \"27-10=\"],\n  [\"52+32=\", \"69-61=\"],\n  [\"97-83=\", \"74-71=\"],\n];\n\nconst body = context.document.body;\nconst missed = [];\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load('items');\n  await context.sync();\n  if (results.items.length > 0) {\n    results.items[0].insertText(newText, Word.InsertLocation.replace);\n    await context.sync();\n  } else {\n    missed.push(oldText);\n  }\n}\nif (missed.length > 0) {\n  console.log('No match found for: ' + missed.join(', '));\n}", "ps1": "$d = $word.ActiveDocument\n\n# Update the date line (first paragraph)\n$d.Paragraphs.Item(1).Range.Text = '2023-04-01 Saturday'\n\n# Update each arithmetic cell in the table, addressed by (row, col)\n$t = $d.Tables.Item(1)\n$values = @(\n  @('89-71=', '37-32=', '95-12=', '1+83=', '91-72='),\n  @('63-32=', '38+27=', '25-19=', '47+51=', '88+8='),\n  @('31+25=', '23+43=', '0+39=', '67+9=', '68-0='),\n  @('1+50=', '59-41=', '76-47=', '43-22=', '83-10='),\n  @('86-74=', '84-67=', '90-41=', '36+44=', '23+62='),\n  @('46-14=', '78-18=', '23+74=', '81-39=', '25+74='),\n  @('44-6=', '98-68=', '70+15=', '56-20=', '61-28='),\n  @('53+37=', '37-12=', '7+45=', '34+23=', '25+47='),\n  @('25+72=', '22+60=', '47-1=', '80-44=', '31+12='),\n  @('43+55=', '66-25=', '41-15=', '21-12=', '27-21='),\n  @('41+5=', '37-12=', '85-83=', '3+10=', '81-19='),\n  @('59+10=', '95-19=', '79+18=', '11+5=', '6+46='),\n  @('21+13=', '67-56=', '0+65=', '1+90=', '77-15='),\n  @('67-18=', '77-14=', '95-6=', '79-50=', '30+48='),\n  @('60+39=', '62+32=', '97-54=', '56-1=', '60+6='),\n  @('9+40=', '96-90=', '76+3=', '87-33=', '90-73='),\n  @('36-22=', '27-5=', '1+58=', '29+24=', '54-27='),\n  @('43+45=', '91-63=', '76-62=', '62-0=', '96-70='),\n  @('30+30=', '10+62=', '70-39=', '9+10=', '91-58='),\n  @('11+12=', '27-10=', '27-10=', '69-61=', '74-71=')\n)\n\nfor ($r = 1; $r -le 20; $r++) {\n  for ($c = 1; $c -le 5; $c++) {\n    $t.Cell($r, $c).Range.Text = $values[$r-1][$c-1]\n  }\n}\n\nWrite-Output \"done\""}
